$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "RP all" header label to reflect the new LLNL comparisons
$ws.Range("J1").Value = "RP all (with new LLNL comparisons)"

# Update the updated V&V results (sigma M / sigma E) for HGL Temperature Rise, RP all column
$ws.Range("J3").Value = 1.27
$ws.Range("K3").Value = 0.59

# Move the active selection as left by the author when saving
$ws.Range("J2").Select() | Out-Null
